$d = $word.ActiveDocument

# 1) Rewrite "Currently, the metrics..." paragraph -> "The metrics measured..."
#    (leading w:tab run is left untouched; only the text after it is replaced)
$found1 = $d.Content.Find.Execute("Currently, the metrics measured for testing is left up to the individual developers and companies, and isn’t regulated by the FDA. However, the two metrics that seem to be most common are McCabe Cyclomatic Complexity and Lines of Code written.", $true, $false, $false, $false, $false, $true, 1, $false, "The metrics measured for testing as of now in the bio-medical field are chosen by the individual developers and companies, and are not regulated by the FDA. However, two metrics that seem to be a common trend, they are McCabe Cyclomatic Complexity and Lines of Code written.", 2)
Write-Host "Replace1:" $found1

# 2) Rewrite "Cyclomatic complexity is measured..." paragraph
$found2 = $d.Content.Find.Execute("Cyclomatic complexity is measured as a way to tell how many different paths the program can take. This is useful in the bio-medical industry, as there should be as few different paths as possible so that the number of places a bug can occur is reduced. For example, in surgical implant technology, under different extreme conditions the device may function very differently in order to keep the patient healthy. This is as expected, but problems arise in testing if all of these different extreme conditions have their own unique branch that the program can proceed down. If makes testing easier, and the program more robust, if when the implant is in certain conditions, it branches off to do tasks specific to those conditions. But when it has completed these, it re-joins the main program flow. This is why McCabe Cyclomatic Complexity is a metric that is currently used, it helps the developers make sure their code doesn’t have too many different branches, where hard to find bugs may exist.", $true, $false, $false, $false, $false, $true, 1, $false, "Cyclomatic complexity is measured as a way to tell how many different paths the program can take. This is useful in the bio-medical industry, as there should be as few different paths as possible so that the number of places a bug can occur is reduced. For example, in surgical implant technology, the device’s behavior varies under differing conditions of extremity, in order to keep the patient healthy. This is expected, but problems arise in testing if all of these different extreme conditions have their own unique branch that the program can proceed down. It makes testing easier and the program more robust, if when the implant is in certain conditions, it branches off to do tasks specific to those conditions. But when it has completed these, it re-joins the main program flow. This is why McCabe Cyclomatic Complexity is a metric that is currently used. It helps the developers make sure their code doesn’t have an excessive amount of branches, where hard to find bugs may exist.", 2)
Write-Host "Replace2:" $found2

# 3) Rewrite "This metric may be useful..." paragraph
$found3 = $d.Content.Find.Execute("This metric may be useful to track in this domain, mostly for the fact that the people who use the code have a vested interest that the code works as advertised, but also most likely don’t have any coding knowledge. Biologists and Doctors may want to be able to look over the code and see if what is needed is being done, since it is their patients whose health may be determined by the software working correctly. Adequate (or excessive) commenting will help these non-programmers better understand the software, and have more peace of mind about its quality.", $true, $false, $false, $false, $false, $true, 1, $false, "This metric may be useful to track in this domain because people who use the code probably have minimal experience in coding, but have a vested interest in the code working as advertised. Biologists and Doctors may want to be able to look over the code and see if what is needed is being done, since it is their patients whose health may be determined by the software working correctly. Adequate (or excessive) commenting will help these non-programmers better understand the software, and have more peace of mind about its quality.", 2)
Write-Host "Replace3:" $found3

# 4) Rewrite "One of the biggest problems..." paragraph.
#    Done as two separate Find/Replace calls that stop right at the
#    w:lastRenderedPageBreak run boundary so that element is preserved,
#    just like in the target revision.
#    (leading w:tab run is left untouched; only the text after it is replaced)
$found4a = $d.Content.Find.Execute("One of the biggest problems in quality assurance in the bio-medical field, in particular surgical equipment and implants, is the inability to test the software under real world conditions. Extensive testing can be done in labs and such, but this can’t guarantee how the software will react under the conditions it will face in the real world. It is extremely difficult to find patients who are willing to have ", $true, $false, $false, $false, $false, $true, 1, $false, "One of the biggest problems in quality assurance in the bio-medical field, in particular surgical equipment and implants, is the inability to test the software under real world conditions. Extensive testing can be done in labs and such, but this can’t guarantee the software will perform correctly under the conditions it will face in the real world. It is extremely difficult to find patients who are willing to ", 2)
Write-Host "Replace4a:" $found4a
$found4b = $d.Content.Find.Execute("surgery performed on them, or implants implanted, for the purpose of testing if the software works correctly or not. Not to mention the massive amount of regulations in regards to testing medical procedures on humans. The software used for these things is tested extensively in labs, under conditions as close to what will happen when used on people as possible, but it is nearly impossible to test under all the conditions that the software may experience in the real world.", $true, $false, $false, $false, $false, $true, 1, $false, "have surgery performed on them, or implants implanted, for the purpose of testing of testing the correctness of software. Not to mention the massive amount of regulations in regards to testing medical procedures on humans. The software used for these purposes is tested extensively in labs, under conditions as close to what will happen when used on humans, but it is nearly impossible to test all of the conditions that the software may be exposed to in the real world.", 2)
Write-Host "Replace4b:" $found4b

# 5) Move the "_GoBack" bookmark from its old location to the empty paragraph
#    right after the "One of the biggest problems..." paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq [char]13) {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text.TrimEnd([char]13) -like "*be exposed to in the real world.") {
            $target = $p
            break
        }
    }
}
if ($target -ne $null) {
    $d.Bookmarks.Add("_GoBack", $target.Range)
    Write-Host "Bookmark added at target paragraph"
} else {
    Write-Host "WARNING: target paragraph for bookmark not found"
}
